# resultados_2023.xlsx - Alpha 1.3: remove stray index column, fix
# known-error columns ("res_c4" / "total") and patch two outlier p39 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The export had an extra, unlabeled pandas-index column in A (styled,
#    values 0..70, no header). Delete it so every other column shifts left
#    by one and column A finally lines up with the real "res_c1" header.
$ws.Columns.Item(1).Delete()

# 2) After the shift, column D is "res_c4" and column E is "total". Every
#    entity had picked up the same (wrong) placeholder of 0 for both; the
#    corrected values are 7.4 and 525.4 across the whole table.
$ws.Range("D2:D72").Value = 7.4
$ws.Range("E2:E72").Value = 525.4

# 3) Column AR is "p39". Two rows had a real (non-zero) score that had been
#    dropped: row 27 (Departamento Administrativo Del Servicio Civil
#    Distrital) and row 72 (Veeduria Distrital).
$ws.Range("AR27").Value = 3.2
$ws.Range("AR72").Value = 4.2
